{"js": "// Update the division-problem answers in the single table of the document.\n// The table has 20 rows x 5 columns; only every 4th row (0-indexed rows\n// 0, 4, 8, 12, 16) actually holds data, the rows in between are blank\n// spacer rows. We address each data cell by its (row, column) position so\n// the edit is unambiguous even though some old/new values repeat elsewhere\n// in the table (e.g. \"31\u00f76=5, 1\" is simultaneously the old text of one cell\n// and the new text of a different cell).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each entry: [rowIndex(0-based), colIndex(0-based), expectedOldText, newText]\nconst updates = [\n  [0, 0, \"87\u00f74=21, 3\", \"98\u00f76=16, 2\"],\n  [0, 1, \"55\u00f75=11, 0\", \"29\u00f73=9, 2\"],\n  [0, 2, \"83\u00f74=20, 3\", \"32\u00f78=4, 0\"],\n  [0, 3, \"85\u00f73=28, 1\", \"41\u00f76=6, 5\"],\n  [0, 4, \"70\u00f78=8, 6\", \"30\u00f74=7, 2\"],\n\n  [4, 0, \"89\u00f78=11, 1\", \"59\u00f78=7, 3\"],\n  [4, 1, \"70\u00f76=11, 4\", \"95\u00f76=15, 5\"],\n  [4, 2, \"74\u00f76=12, 2\", \"26\u00f75=5, 1\"],\n  [4, 3, \"61\u00f78=7, 5\", \"12\u00f77=1, 5\"],\n  [4, 4, \"41\u00f78=5, 1\", \"16\u00f75=3, 1\"],\n\n  [8, 0, \"17\u00f77=2, 3\", \"14\u00f78=1, 6\"],\n  [8, 1, \"31\u00f76=5, 1\", \"69\u00f72=34, 1\"],\n  [8, 2, \"58\u00f76=9, 4\", \"13\u00f72=6, 1\"],\n  [8, 3, \"40\u00f72=20, 0\", \"13\u00f77=1, 6\"],\n  [8, 4, \"57\u00f76=9, 3\", \"34\u00f72=17, 0\"],\n\n  [12, 0, \"68\u00f79=7, 5\", \"48\u00f75=9, 3\"],\n  [12, 1, \"54\u00f76=9, 0\", \"80\u00f79=8, 8\"],\n  [12, 2, \"64\u00f72=32, 0\", \"72\u00f73=24, 0\"],\n  [12, 3, \"98\u00f72=49, 0\", \"39\u00f73=13, 0\"],\n  [12, 4, \"75\u00f74=18, 3\", \"31\u00f76=5, 1\"],\n\n  [16, 0, \"15\u00f77=2, 1\", \"37\u00f74=9, 1\"],\n  [16, 1, \"46\u00f72=23, 0\", \"10\u00f73=3, 1\"],\n  [16, 2, \"54\u00f77=7, 5\", \"19\u00f79=2, 1\"],\n  [16, 3, \"12\u00f78=1, 4\", \"41\u00f79=4, 5\"],\n  [16, 4, \"94\u00f76=15, 4\", \"57\u00f77=8, 1\"],\n];\n\nconst cells = updates.map(([r, c]) => table.getCell(r, c));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  const [, , oldText, newText] = updates[i];\n  const cell = cells[i];\n  // Guard: only overwrite if the current text matches what we expect,\n  // so the script fails loudly instead of silently corrupting the table\n  // if the document layout ever differs from what this script assumes.\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Unexpected cell text at update ${i}: expected \"${oldText}\", found \"${cell.value}\"`\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem answers in the single table of the document.\n# The table has 20 rows x 5 columns; only every 4th row (1-indexed rows\n# 1, 5, 9, 13, 17) actually holds data, the rows in between are blank\n# spacer rows. We address each data cell by its (row, column) position so\n# the edit is unambiguous even though some old/new values repeat elsewhere\n# in the table (e.g. \"31\u00f76=5, 1\" is simultaneously the old text of one cell\n# and the new text of a different cell).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each row: row(1-based), col(1-based), expectedOldText, newText\n$updates = @(\n  @(1, 1, \"87\u00f74=21, 3\", \"98\u00f76=16, 2\"),\n  @(1, 2, \"55\u00f75=11, 0\", \"29\u00f73=9, 2\"),\n  @(1, 3, \"83\u00f74=20, 3\", \"32\u00f78=4, 0\"),\n  @(1, 4, \"85\u00f73=28, 1\", \"41\u00f76=6, 5\"),\n  @(1, 5, \"70\u00f78=8, 6\", \"30\u00f74=7, 2\"),\n\n  @(5, 1, \"89\u00f78=11, 1\", \"59\u00f78=7, 3\"),\n  @(5, 2, \"70\u00f76=11, 4\", \"95\u00f76=15, 5\"),\n  @(5, 3, \"74\u00f76=12, 2\", \"26\u00f75=5, 1\"),\n  @(5, 4, \"61\u00f78=7, 5\", \"12\u00f77=1, 5\"),\n  @(5, 5, \"41\u00f78=5, 1\", \"16\u00f75=3, 1\"),\n\n  @(9, 1, \"17\u00f77=2, 3\", \"14\u00f78=1, 6\"),\n  @(9, 2, \"31\u00f76=5, 1\", \"69\u00f72=34, 1\"),\n  @(9, 3, \"58\u00f76=9, 4\", \"13\u00f72=6, 1\"),\n  @(9, 4, \"40\u00f72=20, 0\", \"13\u00f77=1, 6\"),\n  @(9, 5, \"57\u00f76=9, 3\", \"34\u00f72=17, 0\"),\n\n  @(13, 1, \"68\u00f79=7, 5\", \"48\u00f75=9, 3\"),\n  @(13, 2, \"54\u00f76=9, 0\", \"80\u00f79=8, 8\"),\n  @(13, 3, \"64\u00f72=32, 0\", \"72\u00f73=24, 0\"),\n  @(13, 4, \"98\u00f72=49, 0\", \"39\u00f73=13, 0\"),\n  @(13, 5, \"75\u00f74=18, 3\", \"31\u00f76=5, 1\"),\n\n  @(17, 1, \"15\u00f77=2, 1\", \"37\u00f74=9, 1\"),\n  @(17, 2, \"46\u00f72=23, 0\", \"10\u00f73=3, 1\"),\n  @(17, 3, \"54\u00f77=7, 5\", \"19\u00f79=2, 1\"),\n  @(17, 4, \"12\u00f78=1, 4\", \"41\u00f79=4, 5\"),\n  @(17, 5, \"94\u00f76=15, 4\", \"57\u00f77=8, 1\")\n)\n\nforeach ($u in $updates) {\n  $row = $u[0]\n  $col = $u[1]\n  $oldText = $u[2]\n  $newText = $u[3]\n\n  $cell = $t.Cell($row, $col)\n  # Cell range text includes the trailing cell-mark (\\r\\a); strip it for comparison.\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n  if ($current -ne $oldText) {\n    throw \"Unexpected cell text at row $row col $col`: expected '$oldText', found '$current'\"\n  }\n\n  $cell.Range.Text = $newText\n}\n"}
